$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 2-4 (Bookmark 2p block) and rows 5-6 (Bookmark 2p Packed block) get
# reordered (their text in columns B/D/G/L/M cycles between rows; columns
# A/C/E/F/H/I/J/K are identical across the rows in each block, so they are
# left untouched). We use Copy + PasteSpecial(values) through a scratch area
# far below the used range so that each destination keeps its own original
# cell style/format (PasteSpecial values-only doesn't clobber the
# destination's existing style) while still preserving each source cell's
# original type (text vs number) -- unlike a plain `.Value = ...` assignment,
# which would coerce numeric-looking strings (e.g. "715") into real numbers.

$xlPasteValues = -4163

# Staging rows, well outside the table, used purely as scratch space.
$stage = @{ 2 = 200; 3 = 201; 4 = 202; 5 = 203; 6 = 204 }

foreach ($r in 2..6) {
    $ws.Range("B$r`:M$r").Copy()
    $ws.Range("B$($stage[$r])`:M$($stage[$r])").PasteSpecial($xlPasteValues)
}

# New row 2 = old row 3, new row 3 = old row 4, new row 4 = old row 2
# New row 5 = old row 6, new row 6 = old row 5
$rowMap = @{ 2 = 3; 3 = 4; 4 = 2; 5 = 6; 6 = 5 }

foreach ($destRow in 2..6) {
    $srcRow = $rowMap[$destRow]
    $srcStage = $stage[$srcRow]
    $ws.Range("B$srcStage`:M$srcStage").Copy()
    $ws.Range("B$destRow`:M$destRow").PasteSpecial($xlPasteValues)
}

# Clean up the scratch area so it doesn't linger in the saved workbook.
foreach ($r in $stage.Values) {
    $ws.Range("B$r`:M$r").ClearContents()
}

$excel.CutCopyMode = $false
